$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the date values in J14:J18 while keeping their formatting/style
$ws.Range("J14:J18").Value = $null

# Update the active selection on the sheet to J14:K14
$ws.Activate()
$ws.Range("J14:K14").Select()
